# Update "想去人数" (want-to-go count) figures across the three sheets that
# contain event data. These values were regenerated by the site's data
# pipeline (gh-pages build at 456a3b4), bumping several counters up.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 11980
$ws.Range("F4").Value  = 33
$ws.Range("F5").Value  = 228
$ws.Range("F8").Value  = 11867
$ws.Range("F10").Value = 1177
$ws.Range("F12").Value = 78
$ws.Range("F13").Value = 1790
$ws.Range("F14").Value = 5879
$ws.Range("F15").Value = 129
$ws.Range("F16").Value = 3546
$ws.Range("F17").Value = 194
$ws.Range("F18").Value = 26

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 7

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 11980
$ws.Range("F6").Value  = 33
$ws.Range("F7").Value  = 228
$ws.Range("F8").Value  = 7
$ws.Range("F11").Value = 11867
$ws.Range("F13").Value = 1177
$ws.Range("F15").Value = 78
$ws.Range("F16").Value = 1790
$ws.Range("F18").Value = 5879
$ws.Range("F19").Value = 129
$ws.Range("F20").Value = 3546
$ws.Range("F21").Value = 194
$ws.Range("F22").Value = 26

$wb.Save()
